$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: PESEL -> Pesel
$ws.Range("C1").Value = "Pesel"

# Row 6 (previously Piotr/Kowalski/0430403) becomes Adam/Nowak/<empty>
$ws.Range("A6").Value = "Adam"
$ws.Range("B6").Value = "Nowak"
$ws.Range("C6").Value = ""

# New row 7: Justyna / Justyńska / 4949494 (kept as text, like the other PESEL values)
$ws.Range("A7").Value = "Justyna"
$ws.Range("B7").Value = "Justyńska"
$ws.Range("C7").Value = "'4949494"
